$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1, J1 with same style as existing headers (copy format from H1)
$ws.Cells.Item(1,9).Value = "I0"
$ws.Cells.Item(1,10).Value = "IF"
$ws.Cells.Item(1,8).Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill I and J columns (I0 / IF) for rows 2-65
$ijData = @{
    2 = @(8, 9)
    3 = @(9, 9)
    4 = @(9, 10)
    5 = @(9, 9)
    6 = @(9, 9)
    7 = @(9, 10)
    8 = @(9, 9)
    9 = @(9, 9)
    10 = @(7, 8)
    11 = @(9, 9)
    12 = @(9, 9)
    13 = @(9, 9)
    14 = @(8, 9)
    15 = @(9, 9)
    16 = @(9, 9)
    17 = @(7, 7)
    18 = @(8, 9)
    19 = @(8, 9)
    20 = @(8, 8)
    21 = @(8, 8)
    22 = @(8, 8)
    23 = @(9, 9)
    24 = @(7, 8)
    25 = @(7, 8)
    26 = @(8, 8)
    27 = @(8, 8)
    28 = @(7, 7)
    29 = @(7, 7)
    30 = @(10, 10)
    31 = @(9, 9)
    32 = @(9, 9)
    33 = @(8, 8)
    34 = @(10, 10)
    35 = @(9, 9)
    36 = @(6, 6)
    37 = @(7, 8)
    38 = @(7, 7)
    39 = @(7, 8)
    40 = @(7, 7)
    41 = @(7, 8)
    42 = @(7, 7)
    43 = @(7, 8)
    44 = @(8, 8)
    45 = @(10, 10)
    46 = @(7, 7)
    47 = @(5, 6)
    48 = @(6, 7)
    49 = @(5, 5)
    50 = @(4, 5)
    51 = @(6, 6)
    52 = @(5, 5)
    53 = @(8, 8)
    54 = @(6, 6)
    55 = @(8, 8)
    56 = @(6, 6)
    57 = @(6, 6)
    58 = @(7, 7)
    59 = @(7, 8)
    60 = @(7, 7)
    61 = @(8, 8)
    62 = @(6, 6)
    63 = @(6, 6)
    64 = @(9, 9)
    65 = @(6, 6)
}

foreach ($row in $ijData.Keys) {
    $vals = $ijData[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}

Write-Host "done"